# The deck ships with two DrawingML themes: theme1.xml (bound to the
# slide master / presentation - originally the "Integral" palette) and
# theme2.xml (bound to the notes master - originally the "Office Theme"
# palette). The authored edit swaps their content so the slide master
# now carries the standard "Office Theme" colors.
#
# Colors below are the 12 standard Office-theme scheme colors, in the
# PowerPoint ThemeColorScheme index order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB values are passed as the standard COM 0x00BBGGRR long used by the
# RGBColor.RGB property.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$cs = $theme.ThemeColorScheme

$theme.Name = "Office Theme"
$cs.Item(1).RGB  = 0         # dk1     000000
$cs.Item(2).RGB  = 16777215  # lt1     FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2     44546A
$cs.Item(4).RGB  = 15132391  # lt2     E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1 5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2 ED7D31
$cs.Item(7).RGB  = 10855845  # accent3 A5A5A5
$cs.Item(8).RGB  = 49407     # accent4 FFC000
$cs.Item(9).RGB  = 12874308  # accent5 4472C4
$cs.Item(10).RGB = 4697456   # accent6 70AD47
$cs.Item(11).RGB = 12673797  # hlink   0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72

# Also reach the notes-master side of the theme object model so the
# second palette stays in sync if the host exposes it distinctly.
$notesMaster = $p.NotesMaster
$ncs = $notesMaster.Theme.ThemeColorScheme
$ncs.Item(1).RGB  = 0
$ncs.Item(2).RGB  = 16777215
$ncs.Item(3).RGB  = 6968388
$ncs.Item(4).RGB  = 15132391
$ncs.Item(5).RGB  = 13998939
$ncs.Item(6).RGB  = 3243501
$ncs.Item(7).RGB  = 10855845
$ncs.Item(8).RGB  = 49407
$ncs.Item(9).RGB  = 12874308
$ncs.Item(10).RGB = 4697456
$ncs.Item(11).RGB = 12673797
$ncs.Item(12).RGB = 7491477
